# Insert a new data row before the current row 175, shifting all
# subsequent rows (old 175-286) down to (new 176-287), and populate the
# newly inserted row 175 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 175 (pushes old row175..286 to 176..287)
$ws.Rows("175:175").Insert()

# Populate the new row 175 with the new record
$ws.Cells.Item(175, 1).Value = 4
$ws.Cells.Item(175, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(175, 3).Value = "Los Lagos"
$ws.Cells.Item(175, 4).Value = 44767
$ws.Cells.Item(175, 5).Value = 10
$ws.Cells.Item(175, 6).Value = 100112017
$ws.Cells.Item(175, 7).Value = "Apio"
$ws.Cells.Item(175, 8).Value = "Americana (o)"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 25
$ws.Cells.Item(175, 11).Value = 13000
$ws.Cells.Item(175, 12).Value = 13000
$ws.Cells.Item(175, 13).Value = 13000
$ws.Cells.Item(175, 14).Value = "`$/docena de matas"
$ws.Cells.Item(175, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(175, 16).Value = 2167
$ws.Cells.Item(175, 17).Value = 6
$ws.Cells.Item(175, 18).Value = "Hortaliza"
